$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 97, shifting the existing rows 97-110 down to 98-111.
$ws.Rows.Item(97).EntireRow.Insert()

# Populate the newly-inserted row 97 with the new weekly price record.
$ws.Cells.Item(97, 1).Value = 11
$ws.Cells.Item(97, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(97, 3).Value = "Bíobío"
$ws.Cells.Item(97, 4).Value = 45015
$ws.Cells.Item(97, 5).Value = 8
$ws.Cells.Item(97, 6).Value = 100112012
$ws.Cells.Item(97, 7).Value = "Espinaca"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 150
$ws.Cells.Item(97, 11).Value = 7500
$ws.Cells.Item(97, 12).Value = 8000
$ws.Cells.Item(97, 13).Value = 7667
$ws.Cells.Item(97, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(97, 15).Value = "Región Metropolitana"
$ws.Cells.Item(97, 16).Value = 767
$ws.Cells.Item(97, 17).Value = 10
$ws.Cells.Item(97, 18).Value = "Hortaliza"
